$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 161.279784
$ws.Range("H2").Value = 483.839352
$ws.Range("I2").Value = 0.3023989599621841
$ws.Range("J2").Value = 0.3023989599621841
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 122.328922
$ws.Range("N2").Value = 366.986766
$ws.Range("O2").Value = 0.9783373008518612
$ws.Range("P2").Value = 0.9783373008518613
$ws.Range("Q2").Value = 19729.18211711284
$ws.Range("R2").Value = 177562.6390540156
$ws.Range("S2").Value = 0.2958481822698132
$ws.Range("T2").Value = 0.2958481822698132

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 161.279784
$ws.Range("H3").Value = 483.839352
$ws.Range("I3").Value = 0.3023989599621841
$ws.Range("J3").Value = 0.3023989599621841
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3863573333333334
$ws.Range("N3").Value = 1.159072
$ws.Range("O3").Value = 0.003089929874945324
$ws.Range("P3").Value = 0.003089929874945324
$ws.Range("Q3").Value = 62.311627266816
$ws.Range("R3").Value = 560.804645401344
$ws.Range("S3").Value = 0.0009343915805395473
$ws.Range("T3").Value = 0.0009343915805395473

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 161.279784
$ws.Range("H4").Value = 483.839352
$ws.Range("I4").Value = 0.3023989599621841
$ws.Range("J4").Value = 0.3023989599621841
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.322294
$ws.Range("N4").Value = 6.966882000000001
$ws.Range("O4").Value = 0.0185727692731934
$ws.Range("P4").Value = 0.0185727692731934
$ws.Range("Q4").Value = 374.539074704496
$ws.Range("R4").Value = 3370.851672340464
$ws.Range("S4").Value = 0.005616386111831295
$ws.Range("T4").Value = 0.005616386111831295

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 288.7700093333333
$ws.Range("H5").Value = 866.3100279999999
$ws.Range("I5").Value = 0.541442630470476
$ws.Range("J5").Value = 0.5414426304704759
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 122.328922
$ws.Range("N5").Value = 366.986766
$ws.Range("O5").Value = 0.9783373008518612
$ws.Range("P5").Value = 0.9783373008518613
$ws.Range("Q5").Value = 35324.92394767659
$ws.Range("R5").Value = 317924.3155290894
$ws.Range("S5").Value = 0.5297135216606172
$ws.Range("T5").Value = 0.5297135216606171

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 288.7700093333333
$ws.Range("H6").Value = 866.3100279999999
$ws.Range("I6").Value = 0.541442630470476
$ws.Range("J6").Value = 0.5414426304704759
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3863573333333334
$ws.Range("N6").Value = 1.159072
$ws.Range("O6").Value = 0.003089929874945324
$ws.Range("P6").Value = 0.003089929874945324
$ws.Range("Q6").Value = 111.5684107526684
$ws.Range("R6").Value = 1004.115696774016
$ws.Range("S6").Value = 0.001673019759459705
$ws.Range("T6").Value = 0.001673019759459705

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 288.7700093333333
$ws.Range("H7").Value = 866.3100279999999
$ws.Range("I7").Value = 0.541442630470476
$ws.Range("J7").Value = 0.5414426304704759
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.322294
$ws.Range("N7").Value = 6.966882000000001
$ws.Range("O7").Value = 0.0185727692731934
$ws.Range("P7").Value = 0.0185727692731934
$ws.Range("Q7").Value = 670.608860054744
$ws.Range("R7").Value = 6035.479740492696
$ws.Range("S7").Value = 0.01005608905039907
$ws.Range("T7").Value = 0.01005608905039907

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 83.28466000000002
$ws.Range("H8").Value = 249.85398
$ws.Range("I8").Value = 0.15615840956734
$ws.Range("J8").Value = 0.15615840956734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 122.328922
$ws.Range("N8").Value = 366.986766
$ws.Range("O8").Value = 0.9783373008518612
$ws.Range("P8").Value = 0.9783373008518613
$ws.Range("Q8").Value = 10188.12267693652
$ws.Range("R8").Value = 91693.1040924287
$ws.Range("S8").Value = 0.1527755969214308
$ws.Range("T8").Value = 0.1527755969214308

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 83.28466000000002
$ws.Range("H9").Value = 249.85398
$ws.Range("I9").Value = 0.15615840956734
$ws.Range("J9").Value = 0.15615840956734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3863573333333334
$ws.Range("N9").Value = 1.159072
$ws.Range("O9").Value = 0.003089929874945324
$ws.Range("P9").Value = 0.003089929874945324
$ws.Range("Q9").Value = 32.17763914517334
$ws.Range("R9").Value = 289.5987523065601
$ws.Range("S9").Value = 0.0004825185349460714
$ws.Range("T9").Value = 0.0004825185349460714

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 83.28466000000002
$ws.Range("H10").Value = 249.85398
$ws.Range("I10").Value = 0.15615840956734
$ws.Range("J10").Value = 0.15615840956734
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.322294
$ws.Range("N10").Value = 6.966882000000001
$ws.Range("O10").Value = 0.0185727692731934
$ws.Range("P10").Value = 0.0185727692731934
$ws.Range("Q10").Value = 193.4114662100401
$ws.Range("R10").Value = 1740.703195890361
$ws.Range("S10").Value = 0.002900294110963042
$ws.Range("T10").Value = 0.002900294110963042

